$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 758
$ws.Range("I19").Value = 697.5
$ws.Range("K19").Value = 697.5
$ws.Range("M19").Value = -522.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3539.9666
$ws.Range("I64").Value = 3312.375
$ws.Range("J64").Value = 3622.7273
$ws.Range("K64").Value = 3312.375
$ws.Range("L64").Value = 3622.7273
$ws.Range("M64").Value = -3064.375
$ws.Range("N64").Value = -4118.7273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3539.9666
$ws.Range("I67").Value = 3312.375
$ws.Range("J67").Value = 3622.7273
$ws.Range("K67").Value = 3312.375
$ws.Range("L67").Value = 3622.7273
$ws.Range("M67").Value = -2454.375
$ws.Range("N67").Value = -5338.7273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1594.55
$ws.Range("I137").Value = 1468.3077
$ws.Range("J137").Value = 1829
$ws.Range("K137").Value = 4404.9231
$ws.Range("L137").Value = 5487
$ws.Range("M137").Value = -1854.9231
$ws.Range("N137").Value = -10587

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2700.1755
$ws.Range("I138").Value = 1613.4572
$ws.Range("J138").Value = 4429.0454
$ws.Range("K138").Value = 4840.3716
$ws.Range("L138").Value = 13287.1362
$ws.Range("M138").Value = 299.6283999999996
$ws.Range("N138").Value = -23567.1362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1100.5238
$ws.Range("I74").Value = 943.5625
$ws.Range("J74").Value = 1602.8
$ws.Range("K74").Value = 943.5625
$ws.Range("L74").Value = 1602.8
$ws.Range("M74").Value = -69.5625
$ws.Range("N74").Value = -3350.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1100.5238
$ws.Range("I77").Value = 943.5625
$ws.Range("J77").Value = 1602.8
$ws.Range("K77").Value = 4717.8125
$ws.Range("L77").Value = 8014
$ws.Range("M77").Value = -349.8125
$ws.Range("N77").Value = -16750

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2900
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 2980
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 2980
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -3792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2900
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 2980
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 2980
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -5788

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 41599.5
$ws.Range("J119").Value = 41599.5
$ws.Range("L119").Value = 41599.5
$ws.Range("N119").Value = -51275.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 7346.636
$ws.Range("I122").Value = 10160
$ws.Range("J122").Value = 5002.1665
$ws.Range("K122").Value = 30480
$ws.Range("L122").Value = 15006.4995
$ws.Range("M122").Value = -28030
$ws.Range("N122").Value = -19906.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 24287.166
$ws.Range("J123").Value = 24287.166
$ws.Range("L123").Value = 24287.166
$ws.Range("N123").Value = -34087.166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H131").Value = 33678
$ws.Range("J131").Value = 33678
$ws.Range("L131").Value = 33678
$ws.Range("N131").Value = -43758

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 184538.1
$ws.Range("I86").Value = 3734.1667
$ws.Range("J86").Value = 401502.8
$ws.Range("K86").Value = 3734.1667
$ws.Range("L86").Value = 401502.8
$ws.Range("M86").Value = -2611.1667
$ws.Range("N86").Value = -403748.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 184538.1
$ws.Range("I89").Value = 3734.1667
$ws.Range("J89").Value = 401502.8
$ws.Range("K89").Value = 18670.8335
$ws.Range("L89").Value = 2007514
$ws.Range("M89").Value = -13054.8335
$ws.Range("N89").Value = -2018746

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2894.9
$ws.Range("I31").Value = 1739.0714
$ws.Range("K31").Value = 1739.0714
$ws.Range("M31").Value = -1444.0714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2894.9
$ws.Range("I34").Value = 1739.0714
$ws.Range("K34").Value = 1739.0714
$ws.Range("M34").Value = -1537.0714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 58056.777
$ws.Range("I62").Value = 85534.164
$ws.Range("J62").Value = 3102
$ws.Range("K62").Value = 85534.164
$ws.Range("L62").Value = 3102
$ws.Range("M62").Value = -84910.164
$ws.Range("N62").Value = -4350

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 58056.777
$ws.Range("I65").Value = 85534.164
$ws.Range("J65").Value = 3102
$ws.Range("K65").Value = 427670.82
$ws.Range("L65").Value = 15510
$ws.Range("M65").Value = -424550.82
$ws.Range("N65").Value = -21750

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 31750
$ws.Range("J68").Value = 31750
$ws.Range("L68").Value = 31750
$ws.Range("N68").Value = -33248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 31750
$ws.Range("J71").Value = 31750
$ws.Range("L71").Value = 95250
$ws.Range("N71").Value = -102738

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1608.2683
$ws.Range("I132").Value = 1174.3235
$ws.Range("J132").Value = 3716
$ws.Range("K132").Value = 3522.9705
$ws.Range("L132").Value = 11148
$ws.Range("M132").Value = -992.9704999999999
$ws.Range("N132").Value = -16208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 60000
$ws.Range("J37").Value = 60000
$ws.Range("L37").Value = 180000
$ws.Range("N37").Value = -180224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 586.55554
$ws.Range("I113").Value = 776.5
$ws.Range("K113").Value = 2329.5
$ws.Range("M113").Value = -159.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3488.6965
$ws.Range("I131").Value = 12900
$ws.Range("J131").Value = 1920.1459
$ws.Range("K131").Value = 38700
$ws.Range("L131").Value = 5760.4377
$ws.Range("M131").Value = -33660
$ws.Range("N131").Value = -15840.4377

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 9114.947
$ws.Range("J109").Value = 9114.947
$ws.Range("L109").Value = 9114.947
$ws.Range("N109").Value = -11194.947

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1774.0834
$ws.Range("J113").Value = 2927.8
$ws.Range("L113").Value = 2927.8
$ws.Range("N113").Value = -7267.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2723.8235
$ws.Range("I122").Value = 2351.1667
$ws.Range("J122").Value = 2927.0908
$ws.Range("K122").Value = 7053.500100000001
$ws.Range("L122").Value = 8781.2724
$ws.Range("M122").Value = -4603.500100000001
$ws.Range("N122").Value = -13681.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10893.6
$ws.Range("J123").Value = 10893.6
$ws.Range("L123").Value = 10893.6
$ws.Range("N123").Value = -15793.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 32982.8
$ws.Range("J131").Value = 32982.8
$ws.Range("L131").Value = 32982.8
$ws.Range("N131").Value = -43062.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2552.239
$ws.Range("I132").Value = 2329.9443
$ws.Range("J132").Value = 3352.5
$ws.Range("K132").Value = 6989.8329
$ws.Range("L132").Value = 10057.5
$ws.Range("M132").Value = -4459.8329
$ws.Range("N132").Value = -15117.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2275.25
$ws.Range("I93").Value = 1493
$ws.Range("J93").Value = 2536
$ws.Range("K93").Value = 1493
$ws.Range("L93").Value = 2536
$ws.Range("M93").Value = -245
$ws.Range("N93").Value = -5032

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 22503696
$ws.Range("I122").Value = 19234488
$ws.Range("J122").Value = 28575084
$ws.Range("K122").Value = 57703464
$ws.Range("L122").Value = 85725252
$ws.Range("M122").Value = -57701014
$ws.Range("N122").Value = -85730152

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 19058408
$ws.Range("I136").Value = 25642410
$ws.Range("J136").Value = 717258.2
$ws.Range("K136").Value = 76927230
$ws.Range("L136").Value = 2151774.6
$ws.Range("M136").Value = -76924680
$ws.Range("N136").Value = -2156874.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 134722.62
$ws.Range("J5").Value = 11111.571
$ws.Range("L5").Value = 11111.571
$ws.Range("N5").Value = -11335.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 23000
$ws.Range("J108").Value = 23000
$ws.Range("L108").Value = 23000
$ws.Range("N108").Value = -30680

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 500698
$ws.Range("J119").Value = 500698
$ws.Range("L119").Value = 500698
$ws.Range("N119").Value = -510374

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9767987
$ws.Range("I122").Value = 12501568
$ws.Range("K122").Value = 37504704
$ws.Range("M122").Value = -37502254

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 23913.213
$ws.Range("J123").Value = 23913.213
$ws.Range("L123").Value = 23913.213
$ws.Range("N123").Value = -33713.213

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 54903.57
$ws.Range("J125").Value = 54903.57
$ws.Range("L125").Value = 54903.57
$ws.Range("N125").Value = -64743.57

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1386.52
$ws.Range("I132").Value = 842.7105
$ws.Range("J132").Value = 3108.5833
$ws.Range("K132").Value = 2528.1315
$ws.Range("L132").Value = 9325.749899999999
$ws.Range("M132").Value = 1.86850000000004
$ws.Range("N132").Value = -14385.7499
